# "Generate Report for Handback" - record that the zh-cn and de-de
# handback packages have been generated: populate the "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns
# (J/K/L) on the language sheets, flip the Overview status text, and
# widen the columns that now hold longer values.

$wb = $excel.ActiveWorkbook

$mdTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2e9707f5241826835dedd06eaf74a6d31336198b/e2e/ae14a0ed-0dd7-4940-a671-0137ddd58538.md"
$mdDisplay = "ae14a0ed-0dd7-4940-a671-0137ddd58538.md"

# --- Overview sheet: status text + column widths -----------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("E1").EntireColumn.ColumnWidth = 29.166666666666668
$ws1.Range("F1").EntireColumn.ColumnWidth = 29.166666666666668

# --- zh-cn sheet ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Add($ws2.Range("J2"), $mdTarget, [System.Type]::Missing, [System.Type]::Missing, $mdDisplay)
$ws2.Range("K2").Value = "ae14a0ed-0dd7-4940-a671-0137ddd58538.df8ab9f08ed0f9654e758454c1f61419f0842e70.zh-cn.xlf"
$ws2.Range("L2").Value = "2017-02-17 09:05:37"

$ws2.Hyperlinks.Add($ws2.Range("J3"), $mdTarget, [System.Type]::Missing, [System.Type]::Missing, $mdDisplay)
$ws2.Range("K3").Value = "ae14a0ed-0dd7-4940-a671-0137ddd58538.df8ab9f08ed0f9654e758454c1f61419f0842e70.zh-cn.xlf"
$ws2.Range("L3").Value = "2017-02-17 09:05:37"

$ws2.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$ws2.Range("J1").EntireColumn.ColumnWidth = 39.166666666666664
$ws2.Range("K1").EntireColumn.ColumnWidth = 39.166666666666664

# --- de-de sheet -----------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Add($ws3.Range("J2"), $mdTarget, [System.Type]::Missing, [System.Type]::Missing, $mdDisplay)
$ws3.Range("K2").Value = "ae14a0ed-0dd7-4940-a671-0137ddd58538.df8ab9f08ed0f9654e758454c1f61419f0842e70.de-de.xlf"
$ws3.Range("L2").Value = "2017-02-17 09:06:00"

$ws3.Hyperlinks.Add($ws3.Range("J3"), $mdTarget, [System.Type]::Missing, [System.Type]::Missing, $mdDisplay)
$ws3.Range("K3").Value = "ae14a0ed-0dd7-4940-a671-0137ddd58538.df8ab9f08ed0f9654e758454c1f61419f0842e70.de-de.xlf"
$ws3.Range("L3").Value = "2017-02-17 09:06:00"

$ws3.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$ws3.Range("J1").EntireColumn.ColumnWidth = 39.166666666666664
$ws3.Range("K1").EntireColumn.ColumnWidth = 39.166666666666664

Write-Output "Generate Report for Handback: done"
